# Scheduled-runner market data refresh for the Leve profit workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# on the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets with freshly
# pulled market-board figures, re-deriving the dependent profit values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1250034.5
$ws.Range("I6").Value = 1666679.4
$ws.Range("K6").Value = 5000038.199999999
$ws.Range("M6").Value = -4999926.199999999
$ws.Range("H34").Value = 3253.3333
$ws.Range("I34").Value = 3253.3333
$ws.Range("K34").Value = 3253.3333
$ws.Range("M34").Value = -3050.3333
$ws.Range("H36").Value = 3253.3333
$ws.Range("I36").Value = 3253.3333
$ws.Range("K36").Value = 3253.3333
$ws.Range("M36").Value = -2538.3333
$ws.Range("H40").Value = 5499.75
$ws.Range("J40").Value = 9999
$ws.Range("L40").Value = 9999
$ws.Range("N40").Value = -10349
$ws.Range("H62").Value = 4099.8096
$ws.Range("I62").Value = 4410.9443
$ws.Range("K62").Value = 4410.9443
$ws.Range("M62").Value = -3786.9443
$ws.Range("H65").Value = 4099.8096
$ws.Range("I65").Value = 4410.9443
$ws.Range("K65").Value = 22054.7215
$ws.Range("M65").Value = -18934.7215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 23490.334
$ws.Range("H32").Value = 8783.954
$ws.Range("I32").Value = 5509.525
$ws.Range("K32").Value = 5509.525
$ws.Range("M32").Value = -5222.525
$ws.Range("H45").Value = 25955.75
$ws.Range("I45").Value = 33974.332
$ws.Range("K45").Value = 33974.332
$ws.Range("M45").Value = -33597.332
$ws.Range("H99").Value = 23490.334
$ws.Range("H131").Value = 165987
$ws.Range("J131").Value = 165987
$ws.Range("L131").Value = 165987
$ws.Range("N131").Value = -176067

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 339084.34
$ws.Range("I22").Value = 723.3461
$ws.Range("K22").Value = 723.3461
$ws.Range("M22").Value = -550.3461
$ws.Range("H107").Value = 914.93335
$ws.Range("I107").Value = 809.5769
$ws.Range("J107").Value = 1599.75
$ws.Range("K107").Value = 809.5769
$ws.Range("L107").Value = 1599.75
$ws.Range("M107").Value = 1110.4231
$ws.Range("N107").Value = -5439.75
$ws.Range("H134").Value = 2127.1897
$ws.Range("I134").Value = 2156.74
$ws.Range("J134").Value = 1942.5
$ws.Range("K134").Value = 6470.219999999999
$ws.Range("L134").Value = 5827.5
$ws.Range("M134").Value = -3935.219999999999
$ws.Range("N134").Value = -10897.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 112359.4
$ws.Range("I32").Value = 184666
$ws.Range("J32").Value = 3899.5
$ws.Range("K32").Value = 184666
$ws.Range("L32").Value = 3899.5
$ws.Range("M32").Value = -184350
$ws.Range("N32").Value = -4531.5
$ws.Range("H35").Value = 627.7692
$ws.Range("I35").Value = 627.7692
$ws.Range("K35").Value = 627.7692
$ws.Range("M35").Value = -333.7692
$ws.Range("H42").Value = 4516.3335
$ws.Range("J42").Value = 4516.3335
$ws.Range("L42").Value = 4516.3335
$ws.Range("N42").Value = -5702.3335
$ws.Range("H105").Value = 1282.9546
$ws.Range("I105").Value = 1253.9
$ws.Range("J105").Value = 1307.1666
$ws.Range("K105").Value = 1253.9
$ws.Range("L105").Value = 1307.1666
$ws.Range("M105").Value = 493.0999999999999
$ws.Range("N105").Value = -4801.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 632
$ws.Range("I18").Value = 551.8
$ws.Range("K18").Value = 1655.4
$ws.Range("M18").Value = -1486.4
$ws.Range("H52").Value = 3000
$ws.Range("J52").Value = 3000
$ws.Range("L52").Value = 9000
$ws.Range("N52").Value = -9532
$ws.Range("H61").Value = 150.4375
$ws.Range("I61").Value = 121.5
$ws.Range("J61").Value = 237.25
$ws.Range("K61").Value = 364.5
$ws.Range("L61").Value = 711.75
$ws.Range("M61").Value = -149.5
$ws.Range("N61").Value = -1141.75
$ws.Range("H109").Value = 750
$ws.Range("I109").Value = 750
$ws.Range("K109").Value = 2250
$ws.Range("M109").Value = -1210
$ws.Range("H114").Value = 1426.6666
$ws.Range("I114").Value = 1105.7142
$ws.Range("J114").Value = 2550
$ws.Range("K114").Value = 3317.1426
$ws.Range("L114").Value = 7650
$ws.Range("M114").Value = -63.14259999999967
$ws.Range("N114").Value = -14158
$ws.Range("H121").Value = 672.5
$ws.Range("I121").Value = 30
$ws.Range("K121").Value = 90
$ws.Range("M121").Value = 1220
$ws.Range("H124").Value = 32328.666
$ws.Range("J124").Value = 37230.77
$ws.Range("L124").Value = 111692.31
$ws.Range("N124").Value = -121512.31
$ws.Range("H137").Value = 3384.8572
$ws.Range("I137").Value = 2792.125
$ws.Range("J137").Value = 4175.1665
$ws.Range("K137").Value = 8376.375
$ws.Range("L137").Value = 12525.4995
$ws.Range("M137").Value = -3276.375
$ws.Range("N137").Value = -22725.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 110.76
$ws.Range("J2").Value = 209.33333
$ws.Range("L2").Value = 209.33333
$ws.Range("N2").Value = -435.33333
$ws.Range("H41").Value = 6924.75
$ws.Range("I41").Value = 5850
$ws.Range("J41").Value = 7999.5
$ws.Range("K41").Value = 5850
$ws.Range("L41").Value = 7999.5
$ws.Range("M41").Value = -5495
$ws.Range("N41").Value = -8709.5
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("I102").Value = 1188.7307
$ws.Range("K102").Value = 1188.7307
$ws.Range("M102").Value = 433.2692999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3340.963
$ws.Range("I40").Value = 2473.0454
$ws.Range("K40").Value = 2473.0454
$ws.Range("M40").Value = -2337.0454
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H99").Value = 38000
$ws.Range("J99").Value = 38000
$ws.Range("L99").Value = 38000
$ws.Range("N99").Value = -43990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5974.857
$ws.Range("I2").Value = 5975
$ws.Range("J2").Value = 5974
$ws.Range("K2").Value = 5975
$ws.Range("L2").Value = 5974
$ws.Range("M2").Value = -5863
$ws.Range("N2").Value = -6198
$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20298
$ws.Range("H93").Value = 40389
$ws.Range("J93").Value = 40389
$ws.Range("L93").Value = 40389
$ws.Range("N93").Value = -45381
$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492
$ws.Range("H97").Value = 25814
$ws.Range("J97").Value = 25814
$ws.Range("L97").Value = 25814
$ws.Range("N97").Value = -27796
$ws.Range("H98").Value = 19998.5
$ws.Range("J98").Value = 19998.5
$ws.Range("L98").Value = 19998.5
$ws.Range("N98").Value = -25988.5
$ws.Range("H99").Value = 31249.4
$ws.Range("I99").Value = 25425
$ws.Range("K99").Value = 25425
$ws.Range("M99").Value = -22430
$ws.Range("H100").Value = 617.3333
$ws.Range("I100").Value = 621.3
$ws.Range("J100").Value = 606
$ws.Range("K100").Value = 1242.6
$ws.Range("L100").Value = 1212
$ws.Range("M100").Value = -701.5999999999999
$ws.Range("N100").Value = -2294
